# GPLIM-5973 add aggregation data type to Pooled Tube upload
#
# Adds a new "Aggregation Data Type" column (T) to the PooledTubesTest
# sample sheet, populates the first data row with "Exome", normalizes the
# two "Jon Test Na" sample names to use underscores, and updates the
# window/selection UI state to match the author's saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize sample names (spaces -> underscores)
$ws.Range("B2").Value = "Jon_Test_3a"
$ws.Range("B3").Value = "Jon_Test_4b"

# New column T: "Aggregation Data Type" header + first-row value "Exome"
$ws.Range("T1").Value = "Aggregation Data Type"
$ws.Range("T2").Value = "Exome"

# Give column T an explicit best-fit width like its neighbors
# (COM ColumnWidth is expressed in character units and gets pixel-quantized
# on save; 17.83 is the closest character width to the recorded raw width
# of 18.6640625 twips-per-character units.)
$ws.Columns.Item(20).ColumnWidth = 17.83

# Restore the saved cursor position (cell E1 selected) as in the target file
$ws.Range("E1").Select() | Out-Null

# Match the author's recorded window geometry (best effort; some hosts may
# not persist application-level window metrics back into the saved file)
try {
    $win = $excel.ActiveWindow
    $win.Left = 0
    $win.Top = 0
    $win.Width = 25600
    $win.Height = 16060
} catch {
}
